$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "IT"

$ws.Range("B15").Select()
